$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'1001609528"
$ws.Range("E2").Value = "ACHACALTANASP1"
$ws.Range("H2").Value = "3 jul. 2023, 09:23:51"

$ws.Range("G11").Select()
